$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the template data row (row 2) down into rows 3-10, preserving the
# exact per-cell styles used in row 2 (date style, default text style,
# numeric style) without minting any new style entries.
for ($r = 3; $r -le 10; $r++) {
    $ws.Range("A2:AB2").Copy($ws.Range("A" + $r + ":AB" + $r))
    $ws.Range("A" + $r + ":AB" + $r).RowHeight = 13.5
}

# Rows 11-13 are added as blank but styled rows (same per-column styles,
# no values) - mirror the row-2 formatting, then clear the values out.
for ($r = 11; $r -le 13; $r++) {
    $ws.Range("A2:AB2").Copy($ws.Range("A" + $r + ":AB" + $r))
    $ws.Range("A" + $r + ":AB" + $r).RowHeight = 13.5
    $ws.Range("A" + $r + ":AB" + $r).ClearContents()
}

# Re-create the hyperlinks on the J (email) and K (notification email)
# columns for every newly-populated row, matching the targets already used
# by J2/K2. Hyperlinks.Add reformats the target cell with the built-in
# "Hyperlink" style, so immediately re-copy the plain formatting back from
# the row-2 template cell to keep it looking like J2/K2 (plain text style).
$ws.Hyperlinks.Add($ws.Range("J3"), "mailto:aa@gmail.com")
$ws.Range("J2").Copy($ws.Range("J3"))
$ws.Hyperlinks.Add($ws.Range("K3"), "mailto:aaa@gmail.com")
$ws.Range("K2").Copy($ws.Range("K3"))

$ws.Hyperlinks.Add($ws.Range("J4"), "mailto:aa@gmail.com")
$ws.Range("J2").Copy($ws.Range("J4"))
$ws.Hyperlinks.Add($ws.Range("K4"), "mailto:aaa@gmail.com")
$ws.Range("K2").Copy($ws.Range("K4"))

$ws.Hyperlinks.Add($ws.Range("J5"), "mailto:aa@gmail.com")
$ws.Range("J2").Copy($ws.Range("J5"))
$ws.Hyperlinks.Add($ws.Range("K5"), "mailto:aaa@gmail.com")
$ws.Range("K2").Copy($ws.Range("K5"))

$ws.Hyperlinks.Add($ws.Range("J6"), "mailto:aa@gmail.com")
$ws.Range("J2").Copy($ws.Range("J6"))

$ws.Hyperlinks.Add($ws.Range("J7"), "mailto:aa@gmail.com")
$ws.Range("J2").Copy($ws.Range("J7"))

$ws.Hyperlinks.Add($ws.Range("J8"), "mailto:aa@gmail.com")
$ws.Range("J2").Copy($ws.Range("J8"))

$ws.Hyperlinks.Add($ws.Range("J9"), "mailto:aa@gmail.com")
$ws.Range("J2").Copy($ws.Range("J9"))

$ws.Hyperlinks.Add($ws.Range("J10"), "mailto:aa@gmail.com")
$ws.Range("J2").Copy($ws.Range("J10"))

$ws.Hyperlinks.Add($ws.Range("K6"), "mailto:aaa@gmail.com")
$ws.Range("K2").Copy($ws.Range("K6"))

$ws.Hyperlinks.Add($ws.Range("K7"), "mailto:aaa@gmail.com")
$ws.Range("K2").Copy($ws.Range("K7"))

$ws.Hyperlinks.Add($ws.Range("K8"), "mailto:aaa@gmail.com")
$ws.Range("K2").Copy($ws.Range("K8"))

$ws.Hyperlinks.Add($ws.Range("K9"), "mailto:aaa@gmail.com")
$ws.Range("K2").Copy($ws.Range("K9"))

$ws.Hyperlinks.Add($ws.Range("K10"), "mailto:aaa@gmail.com")
$ws.Range("K2").Copy($ws.Range("K10"))

# Hyperlinks.Add mints a built-in "Hyperlink" named cell style the first
# time it's used; none of the cells actually keep that style (they were
# all reset back to the row-2 formatting above), so drop the now-unused
# named style again.
$wb.Styles.Item("Hyperlink").Delete()

# Match the author's final selection/view state.
$ws.Range("A15").Select()
